$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new text value, taken from the crypto-price refresh diff.
# Some "Price" column values are plain decimals (e.g. "244.27") that Excel
# would otherwise auto-convert to a Number; the source data must stay text
# (it keeps the exact display string, leading zeros, trailing zeros, etc.),
# so those cells get a temporary "@" (Text) number format while the value
# is written, then are restored to the default "Normal" style so no stray
# formatting is left behind.
$updates = @(
    @{ Cell = "D2"; Value = "30.536.70" }
    @{ Cell = "E2"; Value = "  -0.25%  " }
    @{ Cell = "D3"; Value = "1.912.37" }
    @{ Cell = "E3"; Value = "  -0.54%  " }
    @{ Cell = "E4"; Value = "  -0.07%  " }
    @{ Cell = "D5"; Value = "244.27" }
    @{ Cell = "E5"; Value = "  -1.00%  " }
    @{ Cell = "E6"; Value = "  -0.09%  " }
    @{ Cell = "D7"; Value = "0.4853" }
    @{ Cell = "E7"; Value = "  +2.40%  " }
    @{ Cell = "E8"; Value = "  +0.33%  " }
    @{ Cell = "D9"; Value = "0.06796" }
    @{ Cell = "E9"; Value = "  -0.62%  " }
    @{ Cell = "D10"; Value = "111.00" }
    @{ Cell = "E10"; Value = "  +5.63%  " }
    @{ Cell = "D11"; Value = "19.32" }
    @{ Cell = "E11"; Value = "  +5.33%  " }
    @{ Cell = "D12"; Value = "1.911.97" }
    @{ Cell = "E12"; Value = "  -0.52%  " }
    @{ Cell = "D13"; Value = "0.07561" }
    @{ Cell = "E13"; Value = "  -1.73%  " }
    @{ Cell = "D14"; Value = "5.380" }
    @{ Cell = "E14"; Value = "  +0.86%  " }
    @{ Cell = "D15"; Value = "0.6707" }
    @{ Cell = "E15"; Value = "  +0.51%  " }
    @{ Cell = "D16"; Value = "296.50" }
    @{ Cell = "E16"; Value = "  +1.66%  " }
    @{ Cell = "D17"; Value = "30.524.86" }
    @{ Cell = "E17"; Value = "  -0.30%  " }
    @{ Cell = "E18"; Value = "  +0.49%  " }
    @{ Cell = "B19"; Value = "Dai" }
    @{ Cell = "C19"; Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai" }
    @{ Cell = "D19"; Value = "0.9998" }
    @{ Cell = "E19"; Value = "  +0.04%  " }
    @{ Cell = "B20"; Value = "ShibaInu" }
    @{ Cell = "C20"; Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib" }
    @{ Cell = "D20"; Value = "0.000007585" }
    @{ Cell = "E20"; Value = "  -0.47%  " }
    @{ Cell = "D21"; Value = "5.531" }
    @{ Cell = "E21"; Value = "  -0.58%  " }
    @{ Cell = "D22"; Value = "2.161.32" }
    @{ Cell = "E22"; Value = "  -0.36%  " }
    @{ Cell = "D23"; Value = "0.9995" }
    @{ Cell = "E23"; Value = "  -0.17%  " }
    @{ Cell = "D24"; Value = "6.443" }
    @{ Cell = "E24"; Value = "  +0.09%  " }
    @{ Cell = "D25"; Value = "9.465" }
    @{ Cell = "E25"; Value = "  +0.14%  " }
    @{ Cell = "D26"; Value = "165.82" }
    @{ Cell = "E26"; Value = "  -1.03%  " }
    @{ Cell = "E27"; Value = "  -3.32%  " }
    @{ Cell = "D28"; Value = "2.079" }
    @{ Cell = "E28"; Value = "  -1.66%  " }
    @{ Cell = "D29"; Value = "0.1067" }
    @{ Cell = "E29"; Value = "  -0.42%  " }
    @{ Cell = "E30"; Value = "  +2.60%  " }
    @{ Cell = "D31"; Value = "4.147" }
    @{ Cell = "E31"; Value = "  -0.71%  " }
    @{ Cell = "D32"; Value = "4.046" }
    @{ Cell = "E32"; Value = "  -0.34%  " }
    @{ Cell = "D33"; Value = "0.04983" }
    @{ Cell = "E33"; Value = "  -1.14%  " }
    @{ Cell = "D34"; Value = "0.7359" }
    @{ Cell = "E34"; Value = "  -0.21%  " }
    @{ Cell = "D35"; Value = "1.138" }
    @{ Cell = "E35"; Value = "  -0.53%  " }
    @{ Cell = "E36"; Value = "  +0.00%  " }
    @{ Cell = "D37"; Value = "0.02037" }
    @{ Cell = "E37"; Value = "  -1.32%  " }
    @{ Cell = "D38"; Value = "2.715" }
    @{ Cell = "E38"; Value = "  -0.89%  " }
    @{ Cell = "D39"; Value = "2.682" }
    @{ Cell = "E39"; Value = "  -0.21%  " }
    @{ Cell = "E40"; Value = "  -1.70%  " }
    @{ Cell = "D41"; Value = "109.19" }
    @{ Cell = "E41"; Value = "  -1.84%  " }
    @{ Cell = "D42"; Value = "0.4449" }
    @{ Cell = "E42"; Value = "  +1.79%  " }
    @{ Cell = "D43"; Value = "0.8670" }
    @{ Cell = "E43"; Value = "  -0.83%  " }
    @{ Cell = "D44"; Value = "5.801" }
    @{ Cell = "E44"; Value = "  -1.83%  " }
    @{ Cell = "D45"; Value = "0.9994" }
    @{ Cell = "D46"; Value = "69.42" }
    @{ Cell = "E46"; Value = "  +2.12%  " }
    @{ Cell = "D47"; Value = "7.199" }
    @{ Cell = "E47"; Value = "  -1.00%  " }
    @{ Cell = "D48"; Value = "48.40" }
    @{ Cell = "E48"; Value = "  +0.46%  " }
    @{ Cell = "D49"; Value = "9.173" }
    @{ Cell = "E49"; Value = "  -1.62%  " }
    @{ Cell = "D50"; Value = "0.1227" }
    @{ Cell = "E50"; Value = "  -1.18%  " }
    @{ Cell = "D51"; Value = "0.2512" }
    @{ Cell = "E51"; Value = "  -0.55%  " }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $looksNumeric = $u.Value -match "^-?\d+(\.\d+)?$"
    if ($looksNumeric) {
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.Style = "Normal"
    } else {
        $rng.Value = $u.Value
    }
}
